$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitors (Millions)"
$ws.Range("D1").Value = "Year"

# Fill the new Year column for the 20 data rows
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = 2011
}

# Resize column C so its width reflects the new, longer header text
# (closest achievable width to Excel's own best-fit calculation for
# "Overnight International Visitors (Millions)")
$ws.Columns.Item(3).ColumnWidth = 32.3

# Update the sheet selection to match the new state
$ws.Range("D2:D21").Select()
